# Update the "Förändrad" (Changed) date column (C) for every data row
# (rows 2-556) from 2023-09-10 (45179) to 2023-09-11 (45180).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C556").Value = 45180
